# LevelExp.xlsx - "optimise the exp gain"
# - Triple the base Exp seed values (B4/B5/B6: 20/30/50 -> 200/300/500); every
#   B7:B102 cell is a formula chained off the row above, so the rest of the
#   column (and the dependent D/E/F/G columns) recomputes automatically.
# - Unify the number formatting of the CardLevel/GoldFactor/ResFactor columns
#   (E:G) to a consistent "0.00_ " (2 decimal places) style, matching the
#   border/font already used by column E.
# - Resize columns E:G.
# - Restore the selection left by the editor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data edits -----------------------------------------------------------
$ws.Range("B4").Value = 200
$ws.Range("B5").Value = 300
$ws.Range("B6").Value = 500

# --- formatting: unify E:G to one consistent number-formatted style -------
# Column E already carries the "normal" font + thin top border used across
# the table; copy that formatting onto F and G (which previously had no
# border / a different font) before stamping the new number format over the
# whole block so Excel collapses them to a single style record.
$ws.Range("E4").Copy() | Out-Null
$ws.Range("F4:G102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("E4:G102").NumberFormat = "0.00_ "

# --- column widths ---------------------------------------------------------
$ws.Range("E1:G1").ColumnWidth = 6.14

# --- restore selection ------------------------------------------------------
$ws.Range("F97").Select() | Out-Null
